# DRAFT: State of my filesystem
#
# - Refresh the "Estimator" row's sample counts (row 2 for now, before sort)
# - Add a new "Only Possible" row with its own sample counts
# - Add a new "Algorithm" column header in A1
# - Sort the data block by the "Expected" column (I), ascending, which moves
#   "Estimator" to the top and pushes "No Heuristic" to the bottom
# - Leave the selection on A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing "Estimator" row (currently row 5) with fresh counts ---
$ws.Range("A5").Value = "Estimator"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1296
$ws.Range("D5").Value = 14293
$ws.Range("E5").Value = 1482
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = $null
$ws.Range("I5").Formula = "=(B5*1+C5*2+D5*3+E5*4+F5*5+G5*6)/SUM(B5:G5)"

# --- Append a new "Only Possible" row (row 6) ---
$ws.Range("A6").Value = "Only Possible"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1680
$ws.Range("D6").Value = 13053
$ws.Range("E6").Value = 2282
$ws.Range("F6").Value = 57
$ws.Range("G6").Value = 7
$ws.Range("I6").Formula = "=(B6*1+C6*2+D6*3+E6*4+F6*5+G6*6)/SUM(B6:G6)"

# --- Label the new first column ---
$ws.Range("A1").Value = "Algorithm"

# --- Sort the data rows by the Expected column (I), ascending ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("I2:I11"))
$sortObj.SetRange($ws.Range("A2:I11"))
$sortObj.Header = 2
$sortObj.Apply()

# --- Match the saved selection ---
$ws.Range("A2").Select()
